# Rename sheets "wt" -> "wt_log2_expression" and "dcin5" -> "dcin5_log2_expression"
$wb = $excel.ActiveWorkbook

$wtSheet = $wb.Worksheets.Item("wt")
$wtSheet.Name = "wt_log2_expression"

$dcin5Sheet = $wb.Worksheets.Item("dcin5")
$dcin5Sheet.Name = "dcin5_log2_expression"

# Make wt_log2_expression the active/selected tab and update its selected cell
$wtSheet.Select()
$wtSheet.Range("C38").Select()
